$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7645142078399658
$ws.Range("B1").Value = 1.503905415534973
$ws.Range("C1").Value = 5.174703598022461
$ws.Range("D1").Value = 2.9571852684021
$ws.Range("E1").Value = 1.477841258049011
